$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the "Tabela1" table with a new row at the bottom (row 92) via the
# ListObject API, mirroring Excel's native "add table row" behaviour.
$lo = $ws.ListObjects.Item(1)
$newRow = $lo.ListRows.Add()

# Copy the formatting (number format, banded fill, borders, font, alignment)
# from the previous banded data row (row 88) so the new row picks up the
# same "striped" table-row look Excel applies automatically.
$ws.Range("A88:J88").Copy()
$ws.Range("A92:J92").PasteSpecial(-4122)

# New day's Covid-19 data (2020-06-11 update).
$ws.Range("A92").Value2 = 43992
$ws.Range("B92").Value2 = 85626
$ws.Range("C92").Value2 = 758
$ws.Range("D92").Value2 = 1488
$ws.Range("E92").Value2 = 0
$ws.Range("F92").Value2 = 6
$ws.Range("G92").Value2 = 0
$ws.Range("H92").Value2 = 0
$ws.Range("I92").Value2 = 109
$ws.Range("J92").Value2 = 0

# Match the selection recorded in the saved workbook.
$ws.Range("A92:J92").Select()
